# edit.ps1 - PowerPoint COM-interop script
#
# Reproduces the authored change:
#   1. Slide 5's table switches to a different built-in PowerPoint table
#      style (table style gallery pick).
#   2. The deck's theme colour scheme is changed from the "Integral /
#      Red Violet" palette to the built-in "Office" palette (the colours
#      that make up the default "Office Theme").
#
# Helper: build the little-endian BGR integer that PowerPoint's OLE
# automation RGB() values use from a standard RRGGBB hex triple.
function RGBFromHex([int]$r, [int]$g, [int]$b) {
    return ($b * 65536) + ($g * 256) + $r
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Slide 5 - table: apply the new built-in table style.
# ---------------------------------------------------------------------
$slide5 = $p.Slides.Item(5)
for ($i = 1; $i -le $slide5.Shapes.Count; $i++) {
    $shp = $slide5.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{8E1D8935-B36D-40B6-97C3-27C9F46F5558}")
    }
}

# ---------------------------------------------------------------------
# 2) Theme colours - swap the "Integral" palette for the "Office" one.
# ---------------------------------------------------------------------
$officeColors = @(
    @(0,0,0),        # 1  dk1       #000000
    @(255,255,255),  # 2  lt1       #FFFFFF
    @(68,84,106),    # 3  dk2       #44546A
    @(231,230,230),  # 4  lt2       #E7E6E6
    @(91,155,213),   # 5  accent1   #5B9BD5
    @(237,125,49),   # 6  accent2   #ED7D31
    @(165,165,165),  # 7  accent3   #A5A5A5
    @(255,192,0),    # 8  accent4   #FFC000
    @(68,114,196),   # 9  accent5   #4472C4
    @(112,173,71),   # 10 accent6   #70AD47
    @(5,99,193),     # 11 hlink     #0563C1
    @(149,79,114)    # 12 folHlink  #954F72
)

$firstSlide = $p.Slides.Item(1)
for ($i = 1; $i -le $officeColors.Count; $i++) {
    $rgb = $officeColors[$i - 1]
    $firstSlide.ThemeColorScheme.Item($i).RGB = (RGBFromHex $rgb[0] $rgb[1] $rgb[2])
}
